$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # ANLT
$ws2 = $wb.Worksheets.Item(2)   # ANHDT

# Helper: write a date-looking string into a cell without Excel's
# "smart" literal parser turning it into a real date serial number
# (this happens for day<=12 values such as 01/04/2014..04/04/2014,
# which are ambiguous with US-style mm/dd/yyyy parsing). We briefly
# force Text format, assign the value, then restore the original
# date-header formatting by copying it from an already-correct
# header cell (B1), which keeps the same style index (s="1").
function Set-HeaderDate($range, $value, $templateCell) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $templateCell.Copy()
    $range.PasteSpecial(-4122)  # xlPasteFormats
}

# --- Sheet1 (ANLT): shift header dates forward two weeks ---
$ws1.Range("B1").Value = "24/03/2014"
$ws1.Range("C1").Value = "25/03/2014"
$ws1.Range("D1").Value = "26/03/2014"
$ws1.Range("E1").Value = "27/03/2014"
$ws1.Range("F1").Value = "28/03/2014"
$ws1.Range("G1").Value = "29/03/2014"
$ws1.Range("H1").Value = "30/03/2014"
$ws1.Range("I1").Value = "31/03/2014"
Set-HeaderDate $ws1.Range("J1") "01/04/2014" $ws1.Range("B1")
Set-HeaderDate $ws1.Range("K1") "02/04/2014" $ws1.Range("B1")
Set-HeaderDate $ws1.Range("L1") "03/04/2014" $ws1.Range("B1")
Set-HeaderDate $ws1.Range("M1") "04/04/2014" $ws1.Range("B1")

# Grow table1 to include the new row and add the new log entry
$t1 = $ws1.ListObjects.Item(1)
$t1.Resize($ws1.Range("A1:M9"))
$ws1.Range("A3").Value = "Registor"
$ws1.Range("C3").Value = 4

# --- Sheet2 (ANHDT): shift header dates forward two weeks (N1 unchanged) ---
$ws2.Range("B1").Value = "24/03/2014"
$ws2.Range("C1").Value = "25/03/2014"
$ws2.Range("D1").Value = "26/03/2014"
$ws2.Range("E1").Value = "27/03/2014"
$ws2.Range("F1").Value = "28/03/2014"
$ws2.Range("G1").Value = "29/03/2014"
$ws2.Range("H1").Value = "30/03/2014"
$ws2.Range("I1").Value = "31/03/2014"
Set-HeaderDate $ws2.Range("J1") "01/04/2014" $ws2.Range("B1")
Set-HeaderDate $ws2.Range("K1") "02/04/2014" $ws2.Range("B1")
Set-HeaderDate $ws2.Range("L1") "03/04/2014" $ws2.Range("B1")
Set-HeaderDate $ws2.Range("M1") "04/04/2014" $ws2.Range("B1")

# Add the two new log entries (table2 already spans A1:N9)
$ws2.Range("A3").Value = "File URL Controller "
$ws2.Range("A4").Value = "Chỉnh sửa cấu trúc DB"
$ws2.Range("C4").Value = 4

# --- View state: activate ANLT tab, restore selections on each sheet ---
$ws2.Activate()
$ws2.Range("E17").Select()
$ws1.Activate()
$ws1.Range("F20").Select()
